$d = $word.ActiveDocument

# Perform replacements with MatchWholeWord = $true (3rd argument) so that
# short needles like "H1"/"H2"/"40"/"56" can never accidentally match as a
# substring of a longer token (e.g. "H1" inside an already-replaced "H18").
# The order below also guarantees no such collisions occur even without
# MatchWholeWord, since every needle is still untouched at the time it is
# searched for.

# 1. H1 -> H17 (first merge-field result, "NO." label of the first label)
$d.Content.Find.Execute("H1", $true, $true, $false, $false, $false, $true, 1, $false, "H17", 2)

# 2. DANI KRISTIAN -> HENDRI RUSMAWARDANA (name on the first label)
$d.Content.Find.Execute("DANI KRISTIAN", $true, $false, $false, $false, $false, $true, 1, $false, "HENDRI RUSMAWARDANA", 2)

# 3. H2 -> H18 (second merge-field result, "NO." label of the second label)
$d.Content.Find.Execute("H2", $true, $true, $false, $false, $false, $true, 1, $false, "H18", 2)

# 4. SUWANDI -> GANDIS AGUS S. (name on the second label)
$d.Content.Find.Execute("SUWANDI", $true, $false, $false, $false, $false, $true, 1, $false, "GANDIS AGUS S.", 2)

# 5. 40 -> 43 (NO. SHOES value on the second label)
$d.Content.Find.Execute("40", $true, $true, $false, $false, $false, $true, 1, $false, "43", 2)

# 6. 56 -> 58 (NO. TOPI value on the second label)
$d.Content.Find.Execute("56", $true, $true, $false, $false, $false, $true, 1, $false, "58", 2)
